$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Highlight the "Level Design" paragraph under the "#3 (1 week)"
#    section (numId=5) with cyan, matching both the paragraph mark's
#    run properties and the text run's run properties.
# -----------------------------------------------------------------
$rng = $d.Content
$found = $false
$i = 0
while ($rng.Find.Execute("Level Design", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $i = $i + 1
    if ($i -eq 2) {
        $found = $true
        break
    }
    $rng.Collapse(0)
}
if ($found) {
    $para = $rng.Paragraphs(1)
    $para.Range.Font.HighlightColorIndex = 3   # wdTurquoise -> w:highlight val="cyan"
}

# -----------------------------------------------------------------
# 2) Move the lone "_GoBack" bookmark from right after the
#    "Place collectables" run (under the "#4" section, numId=5) to
#    right after the second "Damage received" run that precedes
#    "Death (explosion)" under the Enemies (Robots) bullet of the
#    "#5" section (numId=7). The bookmark is re-created with
#    Bookmarks.Add using the same name, so it automatically vacates
#    its previous spot (Word only keeps one bookmark per name).
#
#    Collapsed ranges placed exactly at a paragraph's trailing
#    boundary are not reliable for Bookmarks.Add in this runtime, so
#    a scratch character is inserted right after the target run,
#    the bookmark is anchored around that scratch character, and the
#    scratch character is then deleted -- leaving the bookmark
#    collapsed in the correct spot, right after the run and before
#    the paragraph end.
# -----------------------------------------------------------------
$rng2 = $d.Content
$target = $null
$j = 0
while ($rng2.Find.Execute("Damage received", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $j = $j + 1
    if ($j -eq 4) {
        $target = $rng2.Duplicate
        break
    }
    $rng2.Collapse(0)
}
if ($target -ne $null) {
    $target.InsertAfter("X")
    $scratch = $d.Range($target.End - 1, $target.End)
    $d.Bookmarks.Add("_GoBack", $scratch) | Out-Null
    $scratch.Text = ""
}
